# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Provincia de Quillota, fecha 44509)
# at the top of the Femacal de La Calera - Tomate block, pushing the
# existing rows 700-726 down to 702-728.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 700 (format copied from
# the row above, matching how the existing data is already styled).
$ws.Rows.Item(700).Insert()
$ws.Rows.Item(701).Insert()

# ---- New row 700 ----
$ws.Cells.Item(700, 1).Value = 3
$ws.Cells.Item(700, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(700, 3).Value = "Coquimbo"
$ws.Cells.Item(700, 4).Value2 = 44509
$ws.Cells.Item(700, 5).Value = 5
$ws.Cells.Item(700, 6).Value = 100112020
$ws.Cells.Item(700, 7).Value = "Tomate"
$ws.Cells.Item(700, 8).Value = "Larga vida"
$ws.Cells.Item(700, 9).Value = "Primera"
$ws.Cells.Item(700, 10).Value = 500
$ws.Cells.Item(700, 11).Value = 17000
$ws.Cells.Item(700, 12).Value = 18000
$ws.Cells.Item(700, 13).Value = 17500
$ws.Cells.Item(700, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(700, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(700, 16).Value = 972
$ws.Cells.Item(700, 17).Value = 18
$ws.Cells.Item(700, 18).Value = "Hortaliza"

# ---- New row 701 ----
$ws.Cells.Item(701, 1).Value = 3
$ws.Cells.Item(701, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(701, 3).Value = "Coquimbo"
$ws.Cells.Item(701, 4).Value2 = 44509
$ws.Cells.Item(701, 5).Value = 5
$ws.Cells.Item(701, 6).Value = 100112020
$ws.Cells.Item(701, 7).Value = "Tomate"
$ws.Cells.Item(701, 8).Value = "Larga vida"
$ws.Cells.Item(701, 9).Value = "Segunda"
$ws.Cells.Item(701, 10).Value = 260
$ws.Cells.Item(701, 11).Value = 14000
$ws.Cells.Item(701, 12).Value = 14000
$ws.Cells.Item(701, 13).Value = 14000
$ws.Cells.Item(701, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(701, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(701, 16).Value = 778
$ws.Cells.Item(701, 17).Value = 18
$ws.Cells.Item(701, 18).Value = "Hortaliza"
